# Hortaliza, Vega Modelo de Temuco - Achicoria
# Insert a new data row at row 63 (shifts existing rows 63-106 down to 64-107)
# and populate it with a new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(63).Insert()

$ws.Range("A63").Value = 10
$ws.Range("B63").Value = "Vega Modelo de Temuco"
$ws.Range("C63").Value = "La Araucanía"
$ws.Range("D63").Value = 45072
$ws.Range("E63").Value = 9
$ws.Range("F63").Value = 100112010
$ws.Range("G63").Value = "Achicoria"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 100
$ws.Range("K63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = 10000
$ws.Range("N63").Value = "$/caja 18 unidades"
$ws.Range("O63").Value = "Región Metropolitana"
$ws.Range("P63").Value = 556
$ws.Range("Q63").Value = 18
$ws.Range("R63").Value = "Hortaliza"
